# Brand new STM32L053R8/Cube1.16.0/Playground/ project.
#
# Update the "Memory" sheet's ROM(KB) "after" figure from 9.86 to 9.85
# (dependent formulas - E9, F9, C14, E14, F14 - recalculate automatically),
# then make "Memory" the active/selected sheet (was "Folder") with its
# selection resting on D11 (was E22); the previously active "Folder" sheet
# keeps its own prior selection (D6) but is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$memory = $wb.Worksheets.Item("Memory")

# Data edit: ROM(KB) "after" value.
$memory.Range("D9").Value = 9.85

# View/selection state: Memory becomes the active sheet & tab, with D11 selected.
$memory.Activate()
[void]$memory.Range("D11").Select()
